# Auto-generated Excel COM-interop script
# Applies scheduled-runner price/profit updates across 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 45999.75
$ws.Range("I54").Value = 50666.332
$ws.Range("K54").Value = 50666.332
$ws.Range("M54").Value = -50180.332

$ws.Range("H74").Value = 5007.4
$ws.Range("I74").Value = 4919.8335
$ws.Range("J74").Value = 5138.75
$ws.Range("K74").Value = 4919.8335
$ws.Range("L74").Value = 5138.75
$ws.Range("M74").Value = -3983.8335
$ws.Range("N74").Value = -7010.75

$ws.Range("H77").Value = 5007.4
$ws.Range("I77").Value = 4919.8335
$ws.Range("J77").Value = 5138.75
$ws.Range("K77").Value = 24599.1675
$ws.Range("L77").Value = 25693.75
$ws.Range("M77").Value = -19919.1675
$ws.Range("N77").Value = -35053.75

$ws.Range("H100").Value = 3367
$ws.Range("J100").Value = 3467.0833
$ws.Range("L100").Value = 3467.0833
$ws.Range("N100").Value = -4549.0833

$ws.Range("H111").Value = 3904.5715
$ws.Range("I111").Value = 3500
$ws.Range("K111").Value = 10500
$ws.Range("M111").Value = -7433

$ws.Range("H129").Value = 1015.2917
$ws.Range("J129").Value = 1043.6428
$ws.Range("L129").Value = 3130.9284
$ws.Range("N129").Value = -13130.9284

$ws.Range("H137").Value = 1231.7451
$ws.Range("I137").Value = 1010.9737
$ws.Range("J137").Value = 1877.0769
$ws.Range("K137").Value = 3032.9211
$ws.Range("L137").Value = 5631.2307
$ws.Range("M137").Value = -482.9211
$ws.Range("N137").Value = -10731.2307

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1456
$ws.Range("I45").Value = 1365.3334
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 1365.3334
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -988.3334
$ws.Range("N45").Value = -2754

$ws.Range("H97").Value = 677.5
$ws.Range("I97").Value = 403.33334
$ws.Range("J97").Value = 1500
$ws.Range("K97").Value = 403.33334
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = 92.66665999999998
$ws.Range("N97").Value = -2492

$ws.Range("H102").Value = 1641.1666
$ws.Range("I102").Value = 1502.3077
$ws.Range("K102").Value = 1502.3077
$ws.Range("M102").Value = 119.6922999999999

$ws.Range("H122").Value = 1328.3143
$ws.Range("I122").Value = 1073.7407
$ws.Range("J122").Value = 2187.5
$ws.Range("K122").Value = 3221.2221
$ws.Range("L122").Value = 6562.5
$ws.Range("M122").Value = -771.2221
$ws.Range("N122").Value = -11462.5

$ws.Range("H132").Value = 2293.4473
$ws.Range("I132").Value = 1591.3572
$ws.Range("J132").Value = 4259.3
$ws.Range("K132").Value = 4774.071599999999
$ws.Range("L132").Value = 12777.9
$ws.Range("M132").Value = -2244.071599999999
$ws.Range("N132").Value = -17837.9

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2955.0667
$ws.Range("J94").Value = 4132.222
$ws.Range("L94").Value = 4132.222
$ws.Range("N94").Value = -5034.222

$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H107").Value = 38046
$ws.Range("I107").Value = 43902.832
$ws.Range("J107").Value = 2905
$ws.Range("K107").Value = 43902.832
$ws.Range("L107").Value = 2905
$ws.Range("M107").Value = -41982.832
$ws.Range("N107").Value = -6745

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 5300
$ws.Range("I86").Value = 4000
$ws.Range("J86").Value = 5671.4287
$ws.Range("K86").Value = 4000
$ws.Range("L86").Value = 5671.4287
$ws.Range("M86").Value = -2877
$ws.Range("N86").Value = -7917.4287

$ws.Range("H89").Value = 5300
$ws.Range("I89").Value = 4000
$ws.Range("J89").Value = 5671.4287
$ws.Range("K89").Value = 20000
$ws.Range("L89").Value = 28357.1435
$ws.Range("M89").Value = -14384
$ws.Range("N89").Value = -39589.14350000001

$ws.Range("H105").Value = 1151
$ws.Range("I105").Value = 1051.3334
$ws.Range("J105").Value = 1450
$ws.Range("K105").Value = 1051.3334
$ws.Range("L105").Value = 1450
$ws.Range("M105").Value = 695.6666
$ws.Range("N105").Value = -4944

$ws.Range("H107").Value = 403
$ws.Range("I107").Value = 452.33334
$ws.Range("K107").Value = 452.33334
$ws.Range("M107").Value = 1467.66666

$ws.Range("H122").Value = 2462.3872
$ws.Range("I122").Value = 2251
$ws.Range("K122").Value = 6753
$ws.Range("M122").Value = -4303

$ws.Range("H135").Value = 40780
$ws.Range("J135").Value = 40780
$ws.Range("L135").Value = 40780
$ws.Range("N135").Value = -50920

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1642.2858
$ws.Range("I5").Value = 1782.5333
$ws.Range("K5").Value = 5347.5999
$ws.Range("M5").Value = -5235.5999

$ws.Range("H107").Value = 325.65
$ws.Range("I107").Value = 437.41666
$ws.Range("J107").Value = 158
$ws.Range("K107").Value = 1312.24998
$ws.Range("L107").Value = 474
$ws.Range("M107").Value = 607.7500199999999
$ws.Range("N107").Value = -4314

$ws.Range("H135").Value = 1642.2858
$ws.Range("I135").Value = 1782.5333
$ws.Range("K135").Value = 16042.7997
$ws.Range("M135").Value = -13507.7997

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H75").Value = 50000
$ws.Range("J75").Value = 50000
$ws.Range("L75").Value = 50000
$ws.Range("N75").Value = -51748

$ws.Range("H78").Value = 50000
$ws.Range("J78").Value = 50000
$ws.Range("L78").Value = 150000
$ws.Range("N78").Value = -158736

$ws.Range("H92").Value = 25757.75
$ws.Range("J92").Value = 25757.75
$ws.Range("L92").Value = 25757.75
$ws.Range("N92").Value = -29501.75

$ws.Range("H95").Value = 1264461
$ws.Range("J95").Value = 1264461
$ws.Range("L95").Value = 1264461
$ws.Range("N95").Value = -1269953

$ws.Range("H97").Value = 42910
$ws.Range("I97").Value = 42910
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 42910
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -42414
$ws.Range("N97").ClearContents()

$ws.Range("H109").Value = 9999.071
$ws.Range("J109").Value = 9999.071
$ws.Range("L109").Value = 9999.071
$ws.Range("N109").Value = -12079.071

$ws.Range("H122").Value = 1639.25
$ws.Range("I122").Value = 1180.3125
$ws.Range("J122").Value = 3475
$ws.Range("K122").Value = 3540.9375
$ws.Range("L122").Value = 10425
$ws.Range("M122").Value = -1090.9375
$ws.Range("N122").Value = -15325

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 17903
$ws.Range("J17").Value = 17903
$ws.Range("L17").Value = 17903
$ws.Range("N17").Value = -18243

$ws.Range("H24").Value = 8000.5
$ws.Range("J24").Value = 8000.5
$ws.Range("L24").Value = 8000.5
$ws.Range("N24").Value = -8686.5

$ws.Range("H68").Value = 2052.5
$ws.Range("I68").Value = 1711.3334
$ws.Range("J68").Value = 2331.6365
$ws.Range("K68").Value = 1711.3334
$ws.Range("L68").Value = 2331.6365
$ws.Range("M68").Value = -962.3334
$ws.Range("N68").Value = -3829.6365

$ws.Range("H71").Value = 2052.5
$ws.Range("I71").Value = 1711.3334
$ws.Range("J71").Value = 2331.6365
$ws.Range("K71").Value = 8556.666999999999
$ws.Range("L71").Value = 11658.1825
$ws.Range("M71").Value = -4812.666999999999
$ws.Range("N71").Value = -19146.1825

$ws.Range("H93").Value = 1200
$ws.Range("J93").Value = 1400
$ws.Range("L93").Value = 1400
$ws.Range("N93").Value = -3896

$ws.Range("H100").Value = 3895.9167
$ws.Range("I100").Value = 4333.4443
$ws.Range("K100").Value = 4333.4443
$ws.Range("M100").Value = -3792.4443

$ws.Range("H122").Value = 17311938
$ws.Range("I122").Value = 14709565
$ws.Range("J122").Value = 22227532
$ws.Range("K122").Value = 44128695
$ws.Range("L122").Value = 66682596
$ws.Range("M122").Value = -44126245
$ws.Range("N122").Value = -66687496

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 711.5263
$ws.Range("I107").Value = 681.26666
$ws.Range("K107").Value = 2043.79998
$ws.Range("M107").Value = -123.79998

$ws.Range("H122").Value = 10778232
$ws.Range("I122").Value = 13160090
$ws.Range("J122").Value = 6252702
$ws.Range("K122").Value = 39480270
$ws.Range("L122").Value = 18758106
$ws.Range("M122").Value = -39477820
$ws.Range("N122").Value = -18763006

$ws.Range("H125").Value = 60712.383
$ws.Range("J125").Value = 60712.383
$ws.Range("L125").Value = 60712.383
$ws.Range("N125").Value = -70552.383
